$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("NewsPaper")

# --- New "Newspaper Vendor" data -------------------------------------------
# Row 2 (Patuli / ToI)
$ws.Range("C2").Value = "PRASENJIT "
$ws.Range("D2").Value = 2000
$ws.Range("H1").Value = "CONTACT NUMBER"
$ws.Range("H2").Value = 7980181141
$ws.Range("G2").Value = 150

# Row 8 (Nayabad / ToI)
$ws.Range("C8").Value = "GANGADHAR SARKAR"
$ws.Range("D8").Value = 3000
$ws.Range("H8").Value = 7003901402
$ws.Range("G8").Value = 150

# Row 4 (Garia / ToI)
$ws.Range("C4").Value = "GAUTAM PATRA"
$ws.Range("D4").Value = 2000
$ws.Range("H4").Value = 9433501891
$ws.Range("G4").Value = 150

# Day of week column filled last across the three rows
$ws.Range("E2").Value = "SUNDAY"
$ws.Range("E4").Value = "SUNDAY"
$ws.Range("E8").Value = "SUNDAY"

# New header (H1) should look like the rest of the header row (bold, shaded)
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false

# Give the new contact-number column a sensible width (matches a best-fit ~11)
$ws.Columns.Item(8).ColumnWidth = 10.166666666666666

# Update selection on the NewsPaper sheet and make it the active tab
$ws.Range("G9").Select()
$ws.Activate()
